# Auto-generated Excel COM-interop script
# Applies targeted cell value changes across multiple worksheets
# as described by the source diff (value-only edits, no formulas).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5138.857
$ws.Range("I2").Value = 1995.6666
$ws.Range("K2").Value = 1995.6666
$ws.Range("M2").Value = -1882.6666
$ws.Range("H12").Value = 131.11111
$ws.Range("I12").Value = 179.8
$ws.Range("K12").Value = 179.8
$ws.Range("M12").Value = -9.800000000000011
$ws.Range("H17").Value = 51500
$ws.Range("J17").Value = 51500
$ws.Range("L17").Value = 154500
$ws.Range("N17").Value = -154836
$ws.Range("H28").Value = 1168
$ws.Range("I28").Value = 1285.4546
$ws.Range("K28").Value = 1285.4546
$ws.Range("M28").Value = -800.4546
$ws.Range("H33").Value = 1877.4
$ws.Range("I33").Value = 2271.75
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 2271.75
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -2042.75
$ws.Range("N33").Value = -758
$ws.Range("H43").Value = 3997.5
$ws.Range("J43").Value = 3995
$ws.Range("L43").Value = 3995
$ws.Range("N43").Value = -4133
$ws.Range("H92").Value = 472.82352
$ws.Range("I92").Value = 424.6154
$ws.Range("K92").Value = 424.6154
$ws.Range("M92").Value = 823.3846
$ws.Range("H101").Value = 12500844
$ws.Range("J101").Value = 964.4286
$ws.Range("L101").Value = 2893.2858
$ws.Range("N101").Value = -6137.2858
$ws.Range("H106").Value = 2172.6667
$ws.Range("I106").Value = 2172.6667
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2172.6667
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -1541.6667
$ws.Range("H107").Value = 376.65
$ws.Range("I107").Value = 379
$ws.Range("J107").Value = 367.25
$ws.Range("K107").Value = 379
$ws.Range("L107").Value = 367.25
$ws.Range("M107").Value = 1541
$ws.Range("N107").Value = -4207.25
$ws.Range("H112").Value = 1425.8182
$ws.Range("J112").Value = 1464.8889
$ws.Range("L112").Value = 4394.6667
$ws.Range("N112").Value = -6610.6667
$ws.Range("H132").Value = 1409.1765
$ws.Range("I132").Value = 1307.4286
$ws.Range("K132").Value = 3922.2858
$ws.Range("M132").Value = -1392.2858
$ws.Range("H137").Value = 1980.2667
$ws.Range("I137").Value = 1669.3077
$ws.Range("J137").Value = 4001.5
$ws.Range("K137").Value = 5007.9231
$ws.Range("L137").Value = 12004.5
$ws.Range("M137").Value = -2457.9231
$ws.Range("N137").Value = -17104.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 207.72728
$ws.Range("I5").Value = 218.6
$ws.Range("K5").Value = 218.6
$ws.Range("M5").Value = -106.6
$ws.Range("H32").Value = 6968.727
$ws.Range("I32").Value = 6968.727
$ws.Range("K32").Value = 6968.727
$ws.Range("M32").Value = -6681.727
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H80").Value = 150000
$ws.Range("J80").Value = 150000
$ws.Range("L80").Value = 150000
$ws.Range("N80").Value = -151996
$ws.Range("H82").Value = 37666
$ws.Range("J82").Value = 37666
$ws.Range("L82").Value = 37666
$ws.Range("N82").Value = -38388
$ws.Range("H83").Value = 150000
$ws.Range("J83").Value = 150000
$ws.Range("L83").Value = 450000
$ws.Range("N83").Value = -459984
$ws.Range("H85").Value = 37666
$ws.Range("J85").Value = 37666
$ws.Range("L85").Value = 37666
$ws.Range("N85").Value = -40162
$ws.Range("H132").Value = 1946.4333
$ws.Range("I132").Value = 1835.8
$ws.Range("K132").Value = 5507.4
$ws.Range("M132").Value = -2977.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 207.72728
$ws.Range("I4").Value = 218.6
$ws.Range("K4").Value = 218.6
$ws.Range("M4").Value = -103.6
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 219.66667
$ws.Range("I12").Value = 219.66667
$ws.Range("K12").Value = 219.66667
$ws.Range("M12").Value = -49.66667000000001
$ws.Range("H22").Value = 8680.666999999999
$ws.Range("I22").Value = 202.42105
$ws.Range("J22").Value = 23324.908
$ws.Range("K22").Value = 202.42105
$ws.Range("L22").Value = 23324.908
$ws.Range("M22").Value = 147.57895
$ws.Range("N22").Value = -24024.908
$ws.Range("H31").Value = 10759.25
$ws.Range("I31").Value = 8997.5
$ws.Range("J31").Value = 11346.5
$ws.Range("K31").Value = 8997.5
$ws.Range("L31").Value = 11346.5
$ws.Range("M31").Value = -8702.5
$ws.Range("N31").Value = -11936.5
$ws.Range("H34").Value = 10759.25
$ws.Range("I34").Value = 8997.5
$ws.Range("J34").Value = 11346.5
$ws.Range("K34").Value = 8997.5
$ws.Range("L34").Value = 11346.5
$ws.Range("M34").Value = -8795.5
$ws.Range("N34").Value = -11750.5
$ws.Range("H122").Value = 1875
$ws.Range("I122").Value = 1666.6666
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4999.9998
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2549.9998
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 2551.111
$ws.Range("I132").Value = 2457.5
$ws.Range("K132").Value = 7372.5
$ws.Range("M132").Value = -4842.5
$ws.Range("H141").Value = 309456.25
$ws.Range("J141").Value = 309456.25
$ws.Range("L141").Value = 309456.25
$ws.Range("N141").Value = -319816.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 12504740
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 20840634
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 62521902
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -62523524
$ws.Range("H71").Value = 12504740
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 20840634
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 187565706
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -187573818
$ws.Range("H76").Value = 10966
$ws.Range("I76").Value = 1899
$ws.Range("K76").Value = 5697
$ws.Range("M76").Value = -5314
$ws.Range("H79").Value = 10966
$ws.Range("I79").Value = 1899
$ws.Range("K79").Value = 5697
$ws.Range("M79").Value = -4371
$ws.Range("H95").Value = 6233.3335
$ws.Range("H129").Value = 2225.4285
$ws.Range("I129").Value = 929.6667
$ws.Range("J129").Value = 10000
$ws.Range("K129").Value = 2789.0001
$ws.Range("L129").Value = 30000
$ws.Range("M129").Value = 2210.9999
$ws.Range("N129").Value = -40000

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 708
$ws.Range("I6").Value = 962
$ws.Range("K6").Value = 962
$ws.Range("M6").Value = -849
$ws.Range("H16").Value = 708
$ws.Range("I16").Value = 962
$ws.Range("K16").Value = 962
$ws.Range("M16").Value = -712
$ws.Range("H62").Value = 89570
$ws.Range("I62").Value = 89570
$ws.Range("K62").Value = 89570
$ws.Range("M62").Value = -88884
$ws.Range("H65").Value = 89570
$ws.Range("I65").Value = 89570
$ws.Range("K65").Value = 268710
$ws.Range("M65").Value = -265278
$ws.Range("H113").Value = 2186.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 29906.379
$ws.Range("J46").Value = 3410.25
$ws.Range("L46").Value = 3410.25
$ws.Range("N46").Value = -3786.25
$ws.Range("H55").Value = 1530.3334
$ws.Range("I55").Value = 1520
$ws.Range("J55").Value = 1539.375
$ws.Range("K55").Value = 1520
$ws.Range("L55").Value = 1539.375
$ws.Range("M55").Value = -1347
$ws.Range("N55").Value = -1885.375
$ws.Range("H82").Value = 5000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 5000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H100").Value = 3817.0908
$ws.Range("I100").Value = 2197.8
$ws.Range("K100").Value = 2197.8
$ws.Range("M100").Value = -1656.8
$ws.Range("H122").Value = 7128.091
$ws.Range("I122").Value = 7942.6665
$ws.Range("J122").Value = 6150.6
$ws.Range("K122").Value = 23827.9995
$ws.Range("L122").Value = 18451.8
$ws.Range("M122").Value = -21377.9995
$ws.Range("N122").Value = -23351.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J49").Value = 43999
$ws.Range("L49").Value = 43999
$ws.Range("N49").Value = -44459
$ws.Range("H54").Value = 45832.832
$ws.Range("J54").Value = 48749.75
$ws.Range("L54").Value = 48749.75
$ws.Range("N54").Value = -49789.75
$ws.Range("H70").Value = 45499.5
$ws.Range("I70").Value = 45499.5
$ws.Range("K70").Value = 45499.5
$ws.Range("M70").Value = -45184.5
$ws.Range("H73").Value = 45499.5
$ws.Range("I73").Value = 45499.5
$ws.Range("K73").Value = 45499.5
$ws.Range("M73").Value = -44407.5
$ws.Range("H107").Value = 628.8333
$ws.Range("I107").Value = 649.5
$ws.Range("J107").Value = 608.1667
$ws.Range("K107").Value = 1948.5
$ws.Range("L107").Value = 1824.5001
$ws.Range("M107").Value = -28.5
$ws.Range("N107").Value = -5664.5001
$ws.Range("H132").Value = 6156.615
$ws.Range("I132").Value = 6156.615
$ws.Range("K132").Value = 18469.845
$ws.Range("M132").Value = -15939.845
